# Finished second cycle of the second sprint (Sprint 1 Phase 2)
# Add the two new tasks ("Descrever Totalmente Primeita/Segunda Feature")
# and mark a day of progress on rows 11 and 12 (F column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Day "F" (3rd tracked day) gets 1 unit of completed effort on the two
# "Sugerir 2 Features" tasks that were finished.
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 1

# Two new rows of work were added to the backlog: describing each feature
# in full, each estimated at 2 units of effort.
$ws.Range("C13").Value = "Descrever Totalmente Primeita Feature"
$ws.Range("D13").Value = 2

$ws.Range("C14").Value = "Descrever Totalmente Segunda Feature"
$ws.Range("D14").Value = 2

# Leave the selection where the user last clicked while editing.
$ws.Range("D13").Select()

$excel.CalculateFull()

# Make sure the burndown chart's cached series data (it plots D18:K18,
# D19:K19 and D20:K20) picks up the new totals.
try {
    $wb.RefreshAll()
} catch {}
try {
    foreach ($co in $ws.ChartObjects()) {
        $co.Chart.Refresh()
    }
} catch {}
